$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'MSG: None

MSG: The decision about which movie to show on Friday did not result in an agreement.
'
$ws.Range("D2").Value = 'no_decision, '

$ws.Range("C3").Value = 'MSG: None

MSG: The decision process did not yield a consensus on which movie to acquire, so no action will be taken regarding the movie selection for Friday.
'
$ws.Range("D3").Value = 'no_decision, '

$ws.Range("C4").Value = 'MSG: None

MSG: The decision-making committee did not reach a conclusion about the movie to show on Friday, and thus the outcome is recorded as "no decision."
'
$ws.Range("D4").Value = 'no_decision, '

$ws.Range("C5").Value = 'MSG: None

MSG: The decision is recorded, and the rights for "Barbie" will be acquired for the screening on Friday.
'

$ws.Range("C6").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made.
'
$ws.Range("D6").Value = 'no_decision, '

$ws.Range("C7").Value = 'MSG: None

MSG: The decision has been recorded as a no-decision outcome for the movie selection.
'
$ws.Range("D7").Value = 'no_decision, '

$ws.Range("C8").Value = 'MSG: None

MSG: The decision to show "Barbie" has been recorded.
'

$ws.Range("C9").Value = 'MSG: None

MSG: The decision has been recorded as no movie was selected in this meeting.
'
$ws.Range("D9").Value = 'no_decision, '

$ws.Range("C10").Value = 'MSG: None

MSG: The decision process has concluded without selecting a movie for Friday.
'
$ws.Range("D10").Value = 'no_decision, '

$ws.Range("C11").Value = 'MSG: None

MSG: The rights to both movies have been acquired for the showing.
'
$ws.Range("D11").Value = 'both_movies, '

$ws.Range("C12").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie".
'

$ws.Range("C13").Value = 'MSG: None

MSG: The committee has decided to select "Barbie" for the Friday showing.
'

$ws.Range("C14").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday has been recorded as "no decision."
'
$ws.Range("D14").Value = 'no_decision, '

$ws.Range("C15").Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Oppenheimer" have been acquired.
'

$ws.Range("C16").Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been recorded.
'

$ws.Range("C17").Value = 'MSG: None

MSG: The conversation ended without a decision regarding which movie to show on Friday. Therefore, I will proceed by calling the no_decision function.
'
$ws.Range("D17").Value = 'no_decision, '

$ws.Range("C18").Value = 'MSG: None

MSG: The rights to both movies have been acquired successfully.
'
$ws.Range("D18").Value = 'both_movies, '

$ws.Range("C19").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has been recorded as no decision being made.
'
$ws.Range("D19").Value = 'no_decision, '

$ws.Range("C20").Value = 'MSG: None

MSG: The committee did not reach a consensus on which movie to show on Friday, so no decision was made regarding the acquisition of movie rights.
'
$ws.Range("D20").Value = 'no_decision, '

$ws.Range("C21").Value = 'MSG: None

MSG: The decision to show a movie on Friday has not been made, so no selection has been finalized.
'
$ws.Range("D21").Value = 'no_decision, '

$ws.Range("C22").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'

$ws.Range("C23").Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding Friday''s movie.
'
$ws.Range("D23").Value = 'no_decision, '

$ws.Range("C24").Value = 'MSG: None

MSG: The decision to show "Barbie" on Friday has been successfully recorded.
'

$ws.Range("C25").Value = 'MSG: None

MSG: The decision to acquire the rights to "Barbie" has been successfully recorded.
'

$ws.Range("C26").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("D26").Value = 'no_decision, '

$ws.Range("C27").Value = 'MSG: None

MSG: The rights for both movies, "Oppenheimer" and "Barbie," have been acquired successfully.
'

$ws.Range("C28").Value = 'MSG: None

MSG: The decision has been recorded, and "Barbie" will be shown on Friday.
'

$ws.Range("C29").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie showing on Friday.
'

$ws.Range("C30").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie".
'

$ws.Range("C31").Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for "Barbie" to be shown on Friday.
'

$ws.Range("C32").Value = 'MSG: None

MSG: I have successfully recorded the decision to acquire the rights for "Barbie" as the movie to be shown on Friday.
'

$ws.Range("C33").Value = 'MSG: None

MSG: The decision has been made that no movie will be selected at this time.
'

$ws.Range("C34").Value = 'MSG: None

MSG: The committee did not reach a decision regarding which movie to show on Friday.
'
$ws.Range("D34").Value = 'no_decision, '

$ws.Range("C35").Value = 'MSG: None

MSG: The decision has been recorded, and "Oppenheimer" will be the movie shown on Friday.
'

$ws.Range("C36").Value = 'MSG: None

MSG: The decision regarding which movie to show on Friday ended without a clear choice being made.
'
$ws.Range("D36").Value = 'no_decision, '

$ws.Range("C37").Value = 'MSG: None

MSG: The rights to "Barbie" have been acquired for Friday''s showing.
'

$ws.Range("C38").Value = 'MSG: None

MSG: The decision has been recorded, and the movie "Oppenheimer" will be acquired for showing on Friday.
'

$ws.Range("C39").Value = 'MSG: None

MSG: The rights to both movies have been acquired for viewing on Friday.
'
$ws.Range("D39").Value = 'both_movies, '

$ws.Range("C40").Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie to be shown on Friday.
'

$ws.Range("C41").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday has been recorded as a no decision.
'
$ws.Range("D41").Value = 'no_decision, '

$ws.Range("C42").Value = 'MSG: None

MSG: The decision has been recorded as no decision made regarding the movie to be shown on Friday.
'
$ws.Range("D42").Value = 'no_decision, '

$ws.Range("C43").Value = 'MSG: None

MSG: The decision was recorded with no movie selected for Friday.
'
$ws.Range("D43").Value = 'no_decision, '

$ws.Range("C44").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has not been finalized, so I have recorded a no decision.
'
$ws.Range("D44").Value = 'no_decision, '

$ws.Range("C45").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie" as the movie to be shown on Friday.
'

$ws.Range("C46").Value = 'MSG: None

MSG: The decision-making process resulted in no agreement on which movie to show on Friday.
'
$ws.Range("D46").Value = 'no_decision, '

$ws.Range("C47").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("D47").Value = 'no_decision, '

$ws.Range("C48").Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired.
'
$ws.Range("D48").Value = 'both_movies, '

$ws.Range("C49").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'

$ws.Range("C50").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'

$ws.Range("C51").Value = 'MSG: None

MSG: The decision process concluded with no consensus on which movie to show, resulting in no decision being made.
'
$ws.Range("D51").Value = 'no_decision, '

$ws.Range("C52").Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range("D52").Value = 'no_decision, '

$ws.Range("C53").Value = 'MSG: None

MSG: The decision has been recorded indicating that no consensus was reached regarding the movie to be shown on Friday.
'
$ws.Range("D53").Value = 'no_decision, '

$ws.Range("C54").Value = 'MSG: None

MSG: The decision about Friday''s movie could not be made, so I will call the no_decision function.
'
$ws.Range("D54").Value = 'no_decision, '

$ws.Range("C55").Value = 'MSG: None

MSG: The decision has been made, and there was no agreement on which movie to show on Friday.
'
$ws.Range("D55").Value = 'no_decision, '

$ws.Range("C56").Value = 'MSG: None

MSG: The decision about the movie to show on Friday has been recorded as no decision.
'
$ws.Range("D56").Value = 'no_decision, '

$ws.Range("C57").Value = 'MSG: None

MSG: The decision has been recorded as no movie will be shown on Friday.
'

$ws.Range("C58").Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired.
'
$ws.Range("D58").Value = 'both_movies, '

$ws.Range("C59").Value = 'MSG: None

MSG: The rights to both movies have been acquired for the showing on Friday.
'
$ws.Range("D59").Value = 'both_movies, '

$ws.Range("C60").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'

$ws.Range("C61").Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired for the presentation on Friday.
'
$ws.Range("D61").Value = 'both_movies, '

$ws.Range("C62").Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding the movie to be shown on Friday.
'
$ws.Range("D62").Value = 'no_decision, '

$ws.Range("C63").Value = 'MSG: None

MSG: The decision process has concluded without a clear choice for Friday''s movie, resulting in no decision being made.
'
$ws.Range("D63").Value = 'no_decision, '

$ws.Range("C64").Value = 'MSG: None

MSG: The decision has been recorded as no movie selection due to the lack of consensus during the discussion.
'

$ws.Range("C65").Value = 'MSG: None

MSG: The rights to both movies have been successfully acquired.
'
$ws.Range("D65").Value = 'both_movies, '

$ws.Range("C66").Value = 'MSG: None

MSG: The decision has been recorded as no decision about the movie for Friday was made.
'
$ws.Range("D66").Value = 'no_decision, '

$ws.Range("C67").Value = 'MSG: None

MSG: The decision regarding Friday''s movie could not be determined, so no action will be taken.
'
$ws.Range("D67").Value = 'no_decision, '

$ws.Range("C68").Value = 'MSG: None

MSG: The decision to acquire the rights for "Oppenheimer" has been successfully recorded.
'
